$d = $word.ActiveDocument
$styles = $d.Styles

# ---------------------------------------------------------------------------
# List1 - drop the <w:contextualSpacing/> toggle from pPr; everything else
# (basedOn/link/rsid/spacing/fonts) is unchanged. Rebuilt via Delete+Add so
# the boolean paragraph toggle can be fully removed instead of merely
# flipped to w:val="0".
# ---------------------------------------------------------------------------
$styles.Item("List1").Delete()
$list1 = $styles.Add("List1", 1)
$list1.NameLocal = "List 1"
$list1.BaseStyle = "List"
$list1.LinkStyle = "List1Char"
$list1.ParagraphFormat.SpaceAfter = 8
$list1.ParagraphFormat.LineSpacingRule = 0
$list1.ParagraphFormat.LineSpacing = 12
$list1.Font.NameAscii = "Times New Roman"
$list1.Font.NameOther = "Times New Roman"
$list1.Font.NameBi = "Times New Roman"
$list1.Font.Size = 12
$list1.Font.SizeBi = 12
$list1.LanguageID = "en"

# ---------------------------------------------------------------------------
# List6 - no longer based on Heading4; drop suppressAutoHyphens,
# contextualSpacing and outlineLvl; rPr becomes plain Times New Roman 12pt
# (no more explicit Bold=False).
# ---------------------------------------------------------------------------
$styles.Item("List6").Delete()
$list6 = $styles.Add("List6", 1)
$list6.NameLocal = "List 6"
$list6.LinkStyle = "List6Char"
$list6.ParagraphFormat.LineSpacingRule = 0
$list6.ParagraphFormat.LineSpacing = 12.95
$list6.ParagraphFormat.FirstLineIndent = 108
$list6.Font.NameAscii = "Times New Roman"
$list6.Font.NameOther = "Times New Roman"
$list6.Font.NameBi = "Times New Roman"
$list6.Font.Size = 12
$list6.Font.SizeBi = 12
$list6.LanguageID = "en"

# ---------------------------------------------------------------------------
# List7 - no longer based on Heading4; drop suppressAutoHyphens and
# outlineLvl (contextualSpacing stays); rPr becomes Times New Roman 12pt,
# explicitly bold.
# ---------------------------------------------------------------------------
$styles.Item("List7").Delete()
$list7 = $styles.Add("List7", 1)
$list7.NameLocal = "List 7"
$list7.LinkStyle = "List7Char"
$list7.ParagraphFormat.SpaceAfter = 8
$list7.ParagraphFormat.LineSpacingRule = 0
$list7.ParagraphFormat.LineSpacing = 12.95
$list7.ParagraphFormat.FirstLineIndent = 162
$list7.NoSpaceBetweenParagraphsOfSameStyle = $true
$list7.Font.NameAscii = "Times New Roman"
$list7.Font.NameOther = "Times New Roman"
$list7.Font.NameBi = "Times New Roman"
$list7.Font.Bold = $true
$list7.Font.Size = 12
$list7.Font.SizeBi = 12
$list7.LanguageID = "en"

# List7Char - flip Bold from False to True in place (position/basedOn kept).
$styles.Item("List7Char").Font.Bold = $true

# ---------------------------------------------------------------------------
# List8 - no longer based on Heading4; drop suppressAutoHyphens,
# contextualSpacing and outlineLvl; rPr becomes plain Times New Roman 12pt
# (no more explicit Bold=False).
# ---------------------------------------------------------------------------
$styles.Item("List8").Delete()
$list8 = $styles.Add("List8", 1)
$list8.NameLocal = "List 8"
$list8.LinkStyle = "List8Char"
$list8.ParagraphFormat.SpaceAfter = 8
$list8.ParagraphFormat.LineSpacingRule = 0
$list8.ParagraphFormat.LineSpacing = 12.95
$list8.ParagraphFormat.FirstLineIndent = 180
$list8.Font.NameAscii = "Times New Roman"
$list8.Font.NameOther = "Times New Roman"
$list8.Font.NameBi = "Times New Roman"
$list8.Font.Size = 12
$list8.Font.SizeBi = 12
$list8.LanguageID = "en"
